$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new daily row (row 43) to the bottom of the data table.
# Column A holds a date formatted as plain text (e.g. "2025/09/30"), not a
# real Excel date, so force text entry and then drop back to the default
# "Normal" style (no explicit number format) to match the rest of the
# column, which was authored without any per-cell style override.
$ws.Range("A43").NumberFormat = "@"
$ws.Range("A43").Value = "2025/10/01"
$ws.Range("A43").Style = "Normal"

$ws.Range("B43").Value = "水"
$ws.Range("C43").Value = 6
$ws.Range("D43").Value = 11
